$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Late" due column), pushing the
# existing N/O/P columns (Late / heading / Outstanding) one column to the right.
$ws.Columns("N").Insert()

# The newly inserted column should take on the same (manually set) width as the
# column to its left (column M, width 10), but without the "best fit" flag.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the sheet's active selection to match the edited state.
$ws.Range("S11").Select() | Out-Null
